$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Paragraph 2: "{m:for i | Sequence{1, 2, 3}}" (was a `for` field)
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML("<w:p $wns>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t>{m:</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>for </w:t></w:r>" +
    "<w:r><w:t>i</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> | </w:t></w:r>" +
    "<w:r><w:t>Sequence{1, 2, 3}}</w:t></w:r>" +
    "</w:p>")

# Paragraph 3: "{m:('dh' + i + '.gif').asImage().setWidth(100)}" (was the image field)
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML("<w:p $wns>" +
    "<w:r><w:t>{</w:t></w:r>" +
    "<w:r><w:t>m</w:t></w:r>" +
    "<w:r><w:t>:</w:t></w:r>" +
    "<w:r><w:t>('dh' + i + '.gif')</w:t></w:r>" +
    "<w:r><w:t>.asImage()</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>.setWidth(100)</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>}</w:t></w:r>" +
    "</w:p>")

# Paragraph 4: "Some text {m:i}" (was the `i` field)
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertXML("<w:p $wns>" +
    "<w:r><w:t>Some text</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t>{m:i</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
    "<w:bookmarkEnd w:id='0'/>" +
    "<w:r><w:t xml:space='preserve'>}</w:t></w:r>" +
    "</w:p>")

# Paragraph 5: "{m:endfor}" (was the `endfor` field)
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertXML("<w:p $wns>" +
    "<w:r><w:t>{</w:t></w:r>" +
    "<w:r><w:t>m:</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>endfor}</w:t></w:r>" +
    "</w:p>")

Write-Host "Final content:" $d.Content.Text
